$wb = $excel.ActiveWorkbook

# "vast" is sheet1 - fill in the newly-added Model-4 (BART) results row (row 6).
$ws1 = $wb.Worksheets.Item("vast")
$ws1.Range("B6").Value = 0.6683
$ws1.Range("C6").Value = 0.7017
$ws1.Range("D6").Value = 0.6683
$ws1.Range("E6").Value = 0.6462
$ws1.Range("F6").Value = 15933.98
$ws1.Range("G6").Value = 401.48
$ws1.Range("H6").Value = 140013315

# Move the active tab / selection from "SEMEVALTASKA" (sheet5) to "vast" (sheet1),
# and update the selected cell on "vast" to H6.
$ws1.Activate()
$ws1.Range("H6").Select() | Out-Null
